# Prefer H1 to <title> if present
#
# 1. Refresh the cached "datetimeFigureOut" field text (2/1/2022 -> 9/8/2022)
#    on the slide master and every slide layout that carries a Date
#    Placeholder with the old cached value.
# 2. Update the explanatory textbox on slide 1: prepend the new “'h1', ”
#    alternative to the existing “'title'” text (first run only, the rest
#    of the runs/formatting are left untouched) and grow the shape to the
#    width PowerPoint's autofit would have produced for the longer line.

$p = $ppt.ActivePresentation

$oldDate = "2/1/2022"
$newDate = "9/8/2022"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master's own Date Placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout hanging off the master.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: the "TextBox 8" callout that lists the fallback order for the
# generated page title.
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(4)

$tr = $shape.TextFrame.TextRange
$firstRun = $tr.Runs(1, 1)
$firstRun.Text = "In order, " + [char]0x2018 + "h1" + [char]0x2019 + ", " + [char]0x2018 + "title" + [char]0x2019 + " "

# Widen the textbox to match the autofit width PowerPoint computes for the
# longer first line (height/position are unaffected).
$shape.Width = 372.64897637795275
